# Applies a fix to the "Recorded By" column (G) on the active worksheet:
# within each comma-separated list of recorder names, the token that is
# exactly "System" (capital S) is swapped into the first position of the
# list (trading places with whatever item currently occupies position 1).
# Cells that do not contain an exact "System" token, or that already have
# "System" in the first position, are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Row + $usedRange.Rows.Count - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $value = $cell.Value2

    if ($null -eq $value) { continue }
    if ($value -notlike "*System*") { continue }

    $parts = $value -split ", "
    if ($parts.Count -le 1) { continue }

    $systemIndex = -1
    for ($i = 0; $i -lt $parts.Count; $i++) {
        if ($parts[$i].Equals("System")) {
            $systemIndex = $i
            break
        }
    }

    if ($systemIndex -le 0) { continue }

    $temp = $parts[0]
    $parts[0] = $parts[$systemIndex]
    $parts[$systemIndex] = $temp

    $cell.Value2 = [string]::Join(", ", $parts)
}
